$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("F15").Value = "The American Journal of Gastroenterology"
$ws.Range("G15").Value = "https://openalex.org/S66441642"
$ws.Range("H15").Value = "Lippincott Williams & Wilkins"
$ws.Range("I15").Value = "0002-9270"

# V15: TRUE -> FALSE, keep it stored as literal text (matches surrounding
# TRUE/FALSE text cells) rather than as a native boolean.
$ws.Range("V15").Formula = "=""FALSE"""
$ws.Range("V15").Copy()
$ws.Range("V15").PasteSpecial(-4163)

# Row 16
$ws.Range("F16").Value = "The American Journal of Gastroenterology"
$ws.Range("G16").Value = "https://openalex.org/S66441642"
$ws.Range("H16").Value = "Lippincott Williams & Wilkins"
$ws.Range("I16").Value = "0002-9270"

# V16: TRUE -> FALSE, same text-literal treatment as V15.
$ws.Range("V16").Formula = "=""FALSE"""
$ws.Range("V16").Copy()
$ws.Range("V16").PasteSpecial(-4163)

$excel.CutCopyMode = 0
